$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.023.06'
$ws.Range('E2').Value = '  +0.06%  '

$ws.Range('D3').Value = '2.300.56'
$ws.Range('E3').Value = '  +0.25%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '300.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.11%  '

$ws.Range('E6').Value = '  -1.42%  '

$ws.Range('E7').Value = '  +3.57%  '

$ws.Range('E8').Value = '  +0.00%  '

$ws.Range('E9').Value = '  +1.55%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.31'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.49%  '

$ws.Range('E11').Value = '  +0.44%  '

$ws.Range('E12').Value = '  +0.63%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '17.80'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.82%  '

$ws.Range('E14').Value = '  -0.51%  '

$ws.Range('D15').Value = '2.658.07'

$ws.Range('D16').Value = '2.305.56'
$ws.Range('E16').Value = '  -1.66%  '

$ws.Range('D18').Value = '42.919.52'
$ws.Range('E18').Value = '  +0.06%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.96'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.45%  '

$ws.Range('D20').Value = '0.0₃0911'
$ws.Range('E20').Value = '  +1.17%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.13'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.50%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.27'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.75%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.81'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.76%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.18'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.14%  '

$ws.Range('E25').Value = '  -0.10%  '

$ws.Range('E26').Value = '  -0.32%  '

$ws.Range('E27').Value = '  -0.16%  '

$ws.Range('E28').Value = '  +0.19%  '

$ws.Range('E29').Value = '  -12.73%  '

$ws.Range('E30').Value = '  +0.48%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '163.27'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.42%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '33.08'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.14%  '

$ws.Range('E33').Value = '  +0.04%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.13'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.29%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '18.20'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.16%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.77'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.25%  '

$ws.Range('E37').Value = '  +0.39%  '

$ws.Range('E38').Value = '  +1.27%  '

$ws.Range('E39').Value = '  +0.86%  '

$ws.Range('E40').Value = '  -0.28%  '

$ws.Range('E41').Value = '  +1.81%  '

$ws.Range('E42').Value = '  -1.29%  '

$ws.Range('D43').Value = '2.017.71'
$ws.Range('E43').Value = '  +2.28%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0286'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.69%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.24'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.11%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.32'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.42%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '17.54'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.51%  '

$ws.Range('E48').Value = '  -1.72%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.29'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.13%  '

$ws.Range('D50').Value = '2.530.55'
$ws.Range('E50').Value = '  +0.45%  '

$ws.Range('E51').Value = '  -0.54%  '
